# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" right after the "总计" summary sheet
#   (and therefore right before the existing "2021-Q3" sheet), with the
#   fund-holding detail table for the new quarter.
# - Insert a new row into the "总计" summary sheet for "2022-Q4" (1 stock,
#   0.02 billion yuan), pushing the existing "2021-Q3" summary row down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$old   = $wb.Worksheets.Item("2021-Q3")

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a fresh row 2 for 2022-Q4, the old
#    2021-Q3 row slides down to row 3 automatically.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Row 2 (new, B2:D2) should carry no special formatting - just clear
# whatever got inherited from the row-insert so it matches a plain cell.
$total.Range("B2:D2").ClearFormats()

# Cell A2 needs to keep the same "row marker" styling the original A2
# carried (bold, thin border, centered) - copy it straight from A3,
# which now holds that original formatting after the row insert.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.02

$total.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with the fund detail table, inserted right
#    before the existing "2021-Q3" sheet.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Header row + the A2 row-marker cell reuse the same bold/border/centered
# style used on "总计"'s own header row.
$total.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$total.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

$q4.Range("A2").Value = 0

# The fund code / size / position figures are stored as text (matches
# the "011501"-style leading-zero codes and "0.60"-style fixed decimals
# used throughout this workbook) - force text so Excel doesn't coerce
# them into numbers and drop the formatting.
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "009140"

$q4.Range("C2").Value = "永赢竞争力精选混合"

$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.60"
$q4.Range("E2").Value = "94.52"
$q4.Range("F2").Value = "2.93"
$q4.Range("G2").Value = "0.0176"

$q4.Range("H2").Value = 10
